$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: 2021年 - fully populated row
$ws.Range("A7").Value = "2021年"
$ws.Range("B7").Value = 100.2
$ws.Range("C7").Value = 101.7
$ws.Range("D7").Value = 101.9
$ws.Range("E7").Value = 101.3
$ws.Range("F7").Value = 101.9
$ws.Range("G7").Value = 101.4
$ws.Range("H7").Value = 101.3
$ws.Range("I7").Value = 102
$ws.Range("J7").Value = 101.9

# Row 8: 2022年 - only column C populated, rest explicitly blank strings
$ws.Range("A8").Value = "2022年"
$ws.Range("B8").Value = ""
$ws.Range("C8").Value = 101.7
$ws.Range("D8").Value = ""
$ws.Range("E8").Value = ""
$ws.Range("F8").Value = ""
$ws.Range("G8").Value = ""
$ws.Range("H8").Value = ""
$ws.Range("I8").Value = ""
$ws.Range("J8").Value = ""

# Copy style from column A label cells (e.g. A6) to A7/A8 so formatting matches
$ws.Range("A6").Copy()
$ws.Range("A7:A8").PasteSpecial(-4122)  # xlPasteFormats
